# Update misc for hlth raster
# - collapse the duplicate "hospital distance" rows (20 & 21) into a single row
#   that carries both the raw distance field and the new "far" risk factor flag
# - rename the underlying column_name from a *_cont_clst to a *_fctb_clst coding
# - remove the now-redundant duplicate row, shrinking the Table1 data range by one row
# - update the active selection left by the editor

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("covars")

# Rename the column_name text that used to read hlthdist_cont_clst
$ws.Range("A20").Value = "hlthdist_fctb_clst"

# Pick up the highlighted-row formatting (font/fill) that lived on row 21's A cell
# so the renamed row keeps that highlight once row 21 disappears.
$ws.Range("A21").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the previously-blank positivefactor / risk_factor_model flags for this row
$ws.Range("E20").Value = "far"
$ws.Range("G20").Value = "y"

# The old row 21 (hlthst_duration_cont_log_scale_clst) was a duplicate entry for the
# same hospital-distance variable; remove it entirely, shifting rows 22:37 up by one.
$ws.Rows.Item(21).Delete()

# Make sure the Table1 list range tracks the new, smaller extent (A1:G36)
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:G36"))

# Restore the view state left behind by the edit
$ws.Range("E21").Select()
